$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 35274.465
$ws.Range("I40").Value = 49378.047
$ws.Range("J40").Value = 2366.111
$ws.Range("K40").Value = 49378.047
$ws.Range("L40").Value = 2366.111
$ws.Range("M40").Value = -49203.047
$ws.Range("N40").Value = -2716.111

$ws.Range("H41").Value = 1441.2
$ws.Range("I41").Value = 2919.8
$ws.Range("K41").Value = 2919.8
$ws.Range("M41").Value = -2479.8

$ws.Range("H46").Value = 996.37036
$ws.Range("J46").Value = 996.37036
$ws.Range("L46").Value = 2989.11108
$ws.Range("N46").Value = -3227.11108

$ws.Range("H60").Value = 996.37036
$ws.Range("J60").Value = 996.37036
$ws.Range("L60").Value = 2989.11108
$ws.Range("N60").Value = -3957.11108

$ws.Range("H62").Value = 2118.6365
$ws.Range("I62").Value = 1930.5
$ws.Range("K62").Value = 1930.5
$ws.Range("M62").Value = -1306.5

$ws.Range("H65").Value = 2118.6365
$ws.Range("I65").Value = 1930.5
$ws.Range("K65").Value = 9652.5
$ws.Range("M65").Value = -6532.5

$ws.Range("H100").Value = 1736.8334
$ws.Range("I100").Value = 1457.8572
$ws.Range("K100").Value = 1457.8572
$ws.Range("M100").Value = -916.8571999999999


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 71572290
$ws.Range("I110").Value = 91091730
$ws.Range("J110").Value = 1004
$ws.Range("K110").Value = 91091730
$ws.Range("L110").Value = 1004
$ws.Range("M110").Value = -91089685
$ws.Range("N110").Value = -5094

$ws.Range("H132").Value = 2337.35
$ws.Range("I132").Value = 1917.4
$ws.Range("K132").Value = 5752.200000000001
$ws.Range("M132").Value = -3222.200000000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 85632.336
$ws.Range("I20").Value = 93226.17999999999
$ws.Range("K20").Value = 93226.17999999999
$ws.Range("M20").Value = -92979.17999999999

$ws.Range("H80").Value = 1365.5714
$ws.Range("I80").Value = 797.2857
$ws.Range("J80").Value = 1649.7142
$ws.Range("K80").Value = 797.2857
$ws.Range("L80").Value = 1649.7142
$ws.Range("M80").Value = 200.7143
$ws.Range("N80").Value = -3645.7142

$ws.Range("H83").Value = 1365.5714
$ws.Range("I83").Value = 797.2857
$ws.Range("J83").Value = 1649.7142
$ws.Range("K83").Value = 3986.4285
$ws.Range("L83").Value = 8248.571
$ws.Range("M83").Value = 1005.5715
$ws.Range("N83").Value = -18232.571

$ws.Range("H94").Value = 747.2857
$ws.Range("I94").Value = 583.125
$ws.Range("J94").Value = 966.1667
$ws.Range("K94").Value = 583.125
$ws.Range("L94").Value = 966.1667
$ws.Range("M94").Value = -132.125
$ws.Range("N94").Value = -1868.1667

$ws.Range("H99").Value = 1373.909
$ws.Range("I99").Value = 1065.9333
$ws.Range("J99").Value = 2033.8572
$ws.Range("K99").Value = 1065.9333
$ws.Range("L99").Value = 2033.8572
$ws.Range("M99").Value = 432.0667000000001
$ws.Range("N99").Value = -5029.8572


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 115.76471
$ws.Range("I7").Value = 61.454544
$ws.Range("J7").Value = 215.33333
$ws.Range("K7").Value = 61.454544
$ws.Range("L7").Value = 215.33333
$ws.Range("M7").Value = 51.545456
$ws.Range("N7").Value = -441.33333

$ws.Range("H16").Value = 842.8570999999999
$ws.Range("I16").Value = 470.33334
$ws.Range("J16").Value = 1122.25
$ws.Range("K16").Value = 470.33334
$ws.Range("L16").Value = 1122.25
$ws.Range("M16").Value = -183.33334
$ws.Range("N16").Value = -1696.25

$ws.Range("H22").Value = 820.1818
$ws.Range("I22").Value = 446.66666
$ws.Range("J22").Value = 2501
$ws.Range("K22").Value = 446.66666
$ws.Range("L22").Value = 2501
$ws.Range("M22").Value = -96.66665999999998
$ws.Range("N22").Value = -3201

$ws.Range("H68").Value = 14353.385
$ws.Range("J68").Value = 14353.385
$ws.Range("L68").Value = 14353.385
$ws.Range("N68").Value = -15851.385

$ws.Range("H71").Value = 14353.385
$ws.Range("J71").Value = 14353.385
$ws.Range("L71").Value = 43060.155
$ws.Range("N71").Value = -50548.155

$ws.Range("H94").Value = 1469.3334
$ws.Range("J94").Value = 1365.5
$ws.Range("L94").Value = 1365.5
$ws.Range("N94").Value = -2267.5

$ws.Range("H95").Value = 10001
$ws.Range("J95").Value = 10001
$ws.Range("L95").Value = 10001
$ws.Range("N95").Value = -15493

$ws.Range("H97").Value = 21500
$ws.Range("I97").Value = 10000
$ws.Range("J97").Value = 33000
$ws.Range("K97").Value = 10000
$ws.Range("L97").Value = 33000
$ws.Range("M97").Value = -9009
$ws.Range("N97").Value = -34982

$ws.Range("H113").Value = 842.8570999999999
$ws.Range("I113").Value = 470.33334
$ws.Range("J113").Value = 1122.25
$ws.Range("K113").Value = 470.33334
$ws.Range("L113").Value = 1122.25
$ws.Range("M113").Value = 1699.66666
$ws.Range("N113").Value = -5462.25

$ws.Range("H122").Value = 557.8182
$ws.Range("I122").Value = 492.47058
$ws.Range("J122").Value = 780
$ws.Range("K122").Value = 1477.41174
$ws.Range("L122").Value = 2340
$ws.Range("M122").Value = 972.58826
$ws.Range("N122").Value = -7240


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 68533.336
$ws.Range("J9").Value = 2000
$ws.Range("L9").Value = 6000
$ws.Range("N9").Value = -6448

$ws.Range("H51").Value = 2962.875
$ws.Range("I51").Value = 2702
$ws.Range("J51").Value = 3049.8333
$ws.Range("K51").Value = 8106
$ws.Range("L51").Value = 9149.499899999999
$ws.Range("M51").Value = -7646
$ws.Range("N51").Value = -10069.4999

$ws.Range("H68").Value = 1751.5
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 3003
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 9009
$ws.Range("M68").Value = -689
$ws.Range("N68").Value = -10631

$ws.Range("H71").Value = 1751.5
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 3003
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 27027
$ws.Range("M71").Value = -444
$ws.Range("N71").Value = -35139

$ws.Range("H131").Value = 824.36
$ws.Range("I131").Value = 526.36365
$ws.Range("J131").Value = 861.19104
$ws.Range("K131").Value = 1579.09095
$ws.Range("L131").Value = 2583.57312
$ws.Range("M131").Value = 3460.90905
$ws.Range("N131").Value = -12663.57312


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2300.3667
$ws.Range("I132").Value = 1646.0454
$ws.Range("J132").Value = 4099.75
$ws.Range("K132").Value = 4938.1362
$ws.Range("L132").Value = 12299.25
$ws.Range("M132").Value = -2408.1362
$ws.Range("N132").Value = -17359.25


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2642.4546
$ws.Range("I61").Value = 1622
$ws.Range("J61").Value = 3492.8333
$ws.Range("K61").Value = 1622
$ws.Range("L61").Value = 3492.8333
$ws.Range("M61").Value = -1420
$ws.Range("N61").Value = -3896.8333

$ws.Range("H96").Value = 13082.333
$ws.Range("I96").Value = 2000
$ws.Range("K96").Value = 2000
$ws.Range("M96").Value = 746

$ws.Range("H97").Value = 19172
$ws.Range("J97").Value = 19172
$ws.Range("L97").Value = 19172
$ws.Range("N97").Value = -21154

$ws.Range("H113").Value = 2642.4546
$ws.Range("I113").Value = 1622
$ws.Range("J113").Value = 3492.8333
$ws.Range("K113").Value = 1622
$ws.Range("L113").Value = 3492.8333
$ws.Range("M113").Value = 548
$ws.Range("N113").Value = -7832.8333


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H113").Value = 569.4
$ws.Range("I113").Value = 425
$ws.Range("J113").Value = 858.2
$ws.Range("K113").Value = 1275
$ws.Range("L113").Value = 2574.6
$ws.Range("M113").Value = 895
$ws.Range("N113").Value = -6914.6

$ws.Range("H119").Value = 40189.8
$ws.Range("J119").Value = 40189.8
$ws.Range("L119").Value = 40189.8
$ws.Range("N119").Value = -49865.8

$ws.Range("H126").Value = 3196.5
$ws.Range("I126").Value = 3234.8
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 9704.400000000001
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -7234.400000000001
$ws.Range("N126").Value = -13955

